$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "91.194.22"
$ws.Cells.Item(2, 5).Value = "  +3.91%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.140.63"
$ws.Cells.Item(3, 5).Value = "  +2.79%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.22%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'220.98"
$ws.Cells.Item(5, 5).Value = "  +6.82%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'622.89"
$ws.Cells.Item(6, 5).Value = "  +1.01%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +4.19%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "'0.918"
$ws.Cells.Item(8, 5).Value = "  +11.75%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  -0.05%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "3.137.64"
$ws.Cells.Item(10, 5).Value = "  +2.74%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.737"
$ws.Cells.Item(11, 5).Value = "  +25.37%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  +6.60%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  +9.03%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'34.61"
$ws.Cells.Item(14, 5).Value = "  +11.04%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "'5.42"
$ws.Cells.Item(15, 5).Value = "  +3.81%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "91.005.44"
$ws.Cells.Item(16, 5).Value = "  +3.95%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "3.726.61"
$ws.Cells.Item(17, 5).Value = "  +2.84%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "3.149.58"
$ws.Cells.Item(18, 5).Value = "  +2.24%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  +22.73%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "'0.0000228"
$ws.Cells.Item(20, 5).Value = "  +10.92%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'14.20"
$ws.Cells.Item(21, 5).Value = "  +9.66%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'433.40"
$ws.Cells.Item(22, 5).Value = "  +3.89%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  +9.30%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  +7.75%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'6.04"
$ws.Cells.Item(25, 5).Value = "  +13.11%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'12.48"
$ws.Cells.Item(26, 5).Value = "  +10.00%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  +3.43%  "

# Row 28
$ws.Cells.Item(28, 2).Value = "WrappedeETH"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Cells.Item(28, 4).Value = "3.307.07"
$ws.Cells.Item(28, 5).Value = "  +2.23%  "

# Row 29
$ws.Cells.Item(29, 2).Value = "Dai"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(29, 4).Value = "'0.999"
$ws.Cells.Item(29, 5).Value = "  -0.04%  "

# Row 30
$ws.Cells.Item(30, 2).Value = "Cronos"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(30, 4).Value = "'0.170"
$ws.Cells.Item(30, 5).Value = "  +8.71%  "

# Row 31
$ws.Cells.Item(31, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(31, 4).Value = "'8.96"
$ws.Cells.Item(31, 5).Value = "  +12.85%  "

# Row 32
$ws.Cells.Item(32, 2).Value = "dogwifhat"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(32, 4).Value = "'3.95"
$ws.Cells.Item(32, 5).Value = "  +13.87%  "

# Row 33
$ws.Cells.Item(33, 2).Value = "Bittensor"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(33, 4).Value = "'530.67"
$ws.Cells.Item(33, 5).Value = "  +6.35%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "Binance-PegBSC-USD"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Cells.Item(34, 4).Value = "'0.884"
$ws.Cells.Item(34, 5).Value = "  -18.69%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "RenderToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Cells.Item(35, 4).Value = "'7.27"
$ws.Cells.Item(35, 5).Value = "  +12.05%  "

# Row 36
$ws.Cells.Item(36, 2).Value = "Fetch.AI"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(36, 4).Value = "'1.32"
$ws.Cells.Item(36, 5).Value = "  +8.58%  "

# Row 37
$ws.Cells.Item(37, 2).Value = "Kaspa"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(37, 4).Value = "'0.142"
$ws.Cells.Item(37, 5).Value = "  +7.83%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  +4.28%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "EthereumClassic"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(39, 4).Value = "'23.50"
$ws.Cells.Item(39, 5).Value = "  +6.80%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "WhiteBITCoin"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Cells.Item(40, 4).Value = "'22.29"
$ws.Cells.Item(40, 5).Value = "  +0.45%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(41, 4).Value = "'1.00"
$ws.Cells.Item(41, 5).Value = "  -0.13%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "Hedera"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(42, 4).Value = "'0.0781"
$ws.Cells.Item(42, 5).Value = "  +18.30%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "Stellar"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(43, 4).Value = "'0.146"
$ws.Cells.Item(43, 5).Value = "  +8.36%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "USDe"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(44, 4).Value = "'1.00"
$ws.Cells.Item(44, 5).Value = "  +0.01%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "'1.94"
$ws.Cells.Item(45, 5).Value = "  +8.61%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "PolygonEcosystemToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Cells.Item(46, 4).Value = "'0.379"
$ws.Cells.Item(46, 5).Value = "  +6.91%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "Monero"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(47, 4).Value = "'143.87"
$ws.Cells.Item(47, 5).Value = "  -3.31%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "OKB"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(48, 4).Value = "'44.17"
$ws.Cells.Item(48, 5).Value = "  +2.03%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "ImmutableX"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(49, 4).Value = "'1.30"
$ws.Cells.Item(49, 5).Value = "  +12.41%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'0.000267"
$ws.Cells.Item(50, 5).Value = "  +27.61%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "Aave"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(51, 4).Value = "'169.69"
$ws.Cells.Item(51, 5).Value = "  +10.11%  "
